# Update cryptocurrency price/volume data per the Jan 14 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell, forcing text (leading apostrophe) so Excel
# doesn't auto-convert number-looking strings and drop significant trailing zeros.
function Set-TextCell($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

$ws.Range("D2").Value = '42.641.74'
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").Value = '2.520.18'
$ws.Range("E3").Value = '  -2.17%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '304.36'
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("D6").Value = '96.65'
$ws.Range("E6").Value = '  -0.56%  '

$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").Value = '0.539'
$ws.Range("E9").Value = '  -2.02%  '

$ws.Range("D10").Value = '36.39'
$ws.Range("E10").Value = '  -0.66%  '

Set-TextCell "D11" '0.0810'
$ws.Range("E11").Value = '  -0.51%  '

$ws.Range("D12").Value = '7.51'
$ws.Range("E12").Value = '  -2.90%  '

$ws.Range("D13").Value = '0.112'
$ws.Range("E13").Value = '  -1.31%  '

$ws.Range("D14").Value = '2.907.83'
$ws.Range("E14").Value = '  -2.10%  '

$ws.Range("D15").Value = '2.509.03'
$ws.Range("E15").Value = '  -2.84%  '

$ws.Range("E16").Value = '  +4.84%  '

$ws.Range("D17").Value = '0.862'
$ws.Range("E17").Value = '  -2.80%  '

$ws.Range("D18").Value = '42.652.54'
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").Value = '12.97'
$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").Value = '0.0₃0973'
$ws.Range("E20").Value = '  -2.39%  '

$ws.Range("D21").Value = '6.45'

$ws.Range("D22").Value = '71.13'
$ws.Range("E22").Value = '  -1.28%  '

$ws.Range("D23").Value = '251.34'
$ws.Range("E23").Value = '  -1.36%  '

$ws.Range("E24").Value = '  -1.15%  '

$ws.Range("D25").Value = '2.03'
$ws.Range("E25").Value = '  -4.28%  '

$ws.Range("D26").Value = '27.03'
$ws.Range("E26").Value = '  -6.34%  '

$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("D28").Value = '2.33'
$ws.Range("E28").Value = '  +10.64%  '

$ws.Range("D29").Value = '10.28'
$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("D30").Value = '38.15'
$ws.Range("E30").Value = '  +0.77%  '

$ws.Range("D31").Value = '5.96'
$ws.Range("E31").Value = '  -1.29%  '

$ws.Range("D32").Value = '155.57'
$ws.Range("E32").Value = '  +0.34%  '

$ws.Range("E33").Value = '  -3.23%  '

$ws.Range("D34").Value = '0.0788'
$ws.Range("E34").Value = '  -2.05%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '2.63'
$ws.Range("E35").Value = '  -4.39%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '2.06'
$ws.Range("E36").Value = '  -5.04%  '

Set-TextCell "D37" '18.60'
$ws.Range("E37").Value = '  +2.54%  '

$ws.Range("E38").Value = '  +1.91%  '

$ws.Range("D39").Value = '24.23'
$ws.Range("E39").Value = '  +4.55%  '

$ws.Range("E40").Value = '  -0.76%  '

$ws.Range("D41").Value = '3.38'
$ws.Range("E41").Value = '  -1.45%  '

$ws.Range("E42").Value = '  -0.23%  '

$ws.Range("E43").Value = '  -0.90%  '

Set-TextCell "D44" '1.00'
$ws.Range("E44").Value = '  +0.13%  '

Set-TextCell "D45" '0.0300'
$ws.Range("E45").Value = '  -3.44%  '

$ws.Range("D46").Value = '2.032.12'
$ws.Range("E46").Value = '  -2.00%  '

$ws.Range("D47").Value = '84.91'
$ws.Range("E47").Value = '  -0.61%  '

$ws.Range("D48").Value = '8.92'
$ws.Range("E48").Value = '  -3.32%  '

$ws.Range("D49").Value = '2.768.31'
$ws.Range("E49").Value = '  -2.12%  '

$ws.Range("D50").Value = '0.189'
$ws.Range("E50").Value = '  -1.08%  '

$ws.Range("D51").Value = '101.75'
$ws.Range("E51").Value = '  -4.65%  '
